$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table-5.1")

# GBR (column D / C2): second & third top sectors were swapped
$ws.Range("D8").Value = "Social, Finance, Analytics, Advertising"
$ws.Range("D9").Value = "Cleantech / Semiconductors "

# GBR (column D / C2): company with highest investment in the (new) second-best sector
$ws.Range("D14").Value = "celltick-technologies"

# Header row: label the country columns (C1/C2/C3 -> C1(USA)/C2(GBR)/C3(IND))
$ws.Range("C4").Value = "C1(USA)"
$ws.Range("D4").Value = "C2(GBR)"
$ws.Range("E4").Value = "C3(IND)"

# Make "Table-5.1" the selected/active sheet & restore its cell selection
$ws.Activate()
$ws.Range("K21").Select() | Out-Null
